$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (not General) interpretation to stay safe for any
# numeric-looking literal by using PowerShell's literal single-quoted
# strings; values that Excel would still auto-parse as numbers get an
# explicit leading apostrophe, matching how a user would type them in
# the Excel UI to keep them as literal text.

$ws.Range("D2").Value = '28.294.43'
$ws.Range("E2").Value = '  +4.28%  '
$ws.Range("D3").Value = '1.728.62'
$ws.Range("E3").Value = '  +2.93%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '''220.56'
$ws.Range("E5").Value = '  +2.60%  '
$ws.Range("E6").Value = '  +0.83%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '''24.29'
$ws.Range("E8").Value = '  +13.40%  '
$ws.Range("D9").Value = '''0.266'
$ws.Range("E9").Value = '  +3.93%  '
$ws.Range("D10").Value = '''0.0635'
$ws.Range("E10").Value = '  +2.10%  '
$ws.Range("D11").Value = '''0.0899'
$ws.Range("E11").Value = '  +1.21%  '
$ws.Range("D12").Value = '1.972.26'
$ws.Range("E12").Value = '  +2.87%  '
$ws.Range("D13").Value = '1.724.95'
$ws.Range("E13").Value = '  +2.90%  '
$ws.Range("E14").Value = '  +3.21%  '
$ws.Range("D15").Value = '''0.560'
$ws.Range("E15").Value = '  +4.45%  '
$ws.Range("D16").Value = '''67.70'
$ws.Range("E16").Value = '  +2.21%  '
$ws.Range("D17").Value = '28.265.57'
$ws.Range("E17").Value = '  +4.16%  '
$ws.Range("D18").Value = '''243.85'
$ws.Range("E18").Value = '  +1.94%  '
$ws.Range("D19").Value = '''8.03'
$ws.Range("E19").Value = '  -0.58%  '
$ws.Range("E20").Value = '  +2.12%  '
$ws.Range("D21").Value = '''1.00'
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("E22").Value = '  +3.00%  '
$ws.Range("D23").Value = '''9.73'
$ws.Range("E23").Value = '  +2.71%  '
$ws.Range("D24").Value = '''2.12'
$ws.Range("E24").Value = '  +0.68%  '
$ws.Range("D25").Value = '''149.60'
$ws.Range("E25").Value = '  +1.68%  '
$ws.Range("D26").Value = '''7.53'
$ws.Range("E26").Value = '  +3.91%  '
$ws.Range("D27").Value = '''16.73'
$ws.Range("E27").Value = '  +2.40%  '
$ws.Range("E28").Value = '  +1.32%  '
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("D30").Value = '''0.0513'
$ws.Range("E30").Value = '  +3.00%  '
$ws.Range("E31").Value = '  +2.40%  '
$ws.Range("E32").Value = '  +2.23%  '
$ws.Range("D33").Value = '1.510.93'
$ws.Range("E33").Value = '  -3.35%  '
$ws.Range("D35").Value = '''1.67'
$ws.Range("E35").Value = '  -1.02%  '
$ws.Range("E36").Value = '  +3.89%  '
$ws.Range("E37").Value = '  +1.35%  '
$ws.Range("D38").Value = '''2.41'
$ws.Range("E38").Value = '  +0.93%  '
$ws.Range("E39").Value = '  +1.03%  '
$ws.Range("E40").Value = '  +0.99%  '
$ws.Range("D41").Value = '''70.97'
$ws.Range("E41").Value = '  +2.85%  '
$ws.Range("D42").Value = '''5.74'
$ws.Range("E42").Value = '  +2.95%  '
$ws.Range("D43").Value = '''1.00'
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("E44").Value = '  +2.05%  '
$ws.Range("D45").Value = '1.876.47'
$ws.Range("E45").Value = '  +2.78%  '
$ws.Range("D46").Value = '''0.807'
$ws.Range("E46").Value = '  +2.79%  '
$ws.Range("D47").Value = '''1.75'
$ws.Range("E47").Value = '  +10.91%  '
$ws.Range("D48").Value = '''90.99'
$ws.Range("E48").Value = '  +0.40%  '
$ws.Range("D49").Value = '0.0₆0112'
$ws.Range("E49").Value = '  +4.56%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '''0.105'
$ws.Range("E50").Value = '  +1.75%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '''8.24'
$ws.Range("E51").Value = '  +1.02%  '
